$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 111473793
$ws.Range("B3").Value = 93388
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 2180
$ws.Range("F3").Value = "Blåmossa"
$ws.Range("G3").Value = "Leucobryum glaucum"
$ws.Range("H3").Value = "(Hedw.) Ångstr."
$ws.Range("Q3").Value = 703959.3331032015
$ws.Range("R3").Value = 6572805.612961343
$ws.Range("AO3").ClearContents()

# Row 4
$ws.Range("A4").Value = 111473776
$ws.Range("B4").Value = 89405
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 1202
$ws.Range("F4").Value = "Ullticka"
$ws.Range("G4").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H4").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q4").Value = 703970.8884549731
$ws.Range("R4").Value = 6572810.333898042
$ws.Range("AO4").Value = "granlåga"

# Row 5
$ws.Range("A5").Value = 111473777
$ws.Range("B5").Value = 89425
$ws.Range("E5").Value = 5442
$ws.Range("F5").Value = "Tallticka"
$ws.Range("G5").Value = "Porodaedalea pini"
$ws.Range("H5").Value = "(Brot.) Murrill"
$ws.Range("Q5").Value = 704301.1177162804
$ws.Range("R5").Value = 6573209.392206083
$ws.Range("AO5").Value = "gammeltall"

# Row 7
$ws.Range("A7").Value = 111473779
$ws.Range("B7").Value = 89425
$ws.Range("E7").Value = 5442
$ws.Range("F7").Value = "Tallticka"
$ws.Range("G7").Value = "Porodaedalea pini"
$ws.Range("H7").Value = "(Brot.) Murrill"
$ws.Range("Q7").Value = 704193.4830821306
$ws.Range("R7").Value = 6572948.378178579
$ws.Range("AO7").Value = "gammeltall"

# Row 8
$ws.Range("A8").Value = 111473773
$ws.Range("Q8").Value = 704016.0051346947
$ws.Range("R8").Value = 6572801.994589122

# Row 9
$ws.Range("A9").Value = 111473791
$ws.Range("B9").Value = 93289
$ws.Range("E9").Value = 2170
$ws.Range("F9").Value = "Flagellkvastmossa"
$ws.Range("G9").Value = "Dicranum flagellare"
$ws.Range("H9").Value = "Hedw."
$ws.Range("Q9").Value = 704004.9502936595
$ws.Range("R9").Value = 6572835.740028554
$ws.Range("AO9").Value = "låga av tall"

# Row 10
$ws.Range("A10").Value = 111473782
$ws.Range("B10").Value = 89183
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 3215
$ws.Range("F10").Value = "Rödgul trumpetsvamp"
$ws.Range("G10").Value = "Craterellus lutescens"
$ws.Range("H10").Value = "(Fr.) Fr."
$ws.Range("Q10").Value = 704171.5165585374
$ws.Range("R10").Value = 6572850.843097115
$ws.Range("AO10").ClearContents()

# Row 11
$ws.Range("A11").Value = 111473784
$ws.Range("B11").Value = 73634
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 6426
$ws.Range("F11").Value = "Kattfotslav"
$ws.Range("G11").Value = "Felipes leucopellaeus"
$ws.Range("H11").Value = "(Ach.) Frisch & G.Thor"
$ws.Range("Q11").Value = 704135.470341172
$ws.Range("R11").Value = 6572843.267234835
$ws.Range("AO11").Value = "äldre gran"

# Row 12
$ws.Range("A12").Value = 111473774
$ws.Range("B12").Value = 89405
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 1202
$ws.Range("F12").Value = "Ullticka"
$ws.Range("G12").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H12").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("M12").ClearContents()
$ws.Range("Q12").Value = 703999.5190368021
$ws.Range("R12").Value = 6572850.823973293
$ws.Range("AO12").Value = "granlåga"

# Row 13
$ws.Range("A13").Value = 111473775
$ws.Range("B13").Value = 89405
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 1202
$ws.Range("F13").Value = "Ullticka"
$ws.Range("G13").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H13").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q13").Value = 703969.3444121893
$ws.Range("R13").Value = 6572791.287347207
$ws.Range("AO13").Value = "granlåga"

# Row 14
$ws.Range("A14").Value = 111473792
$ws.Range("B14").Value = 5113
$ws.Range("E14").Value = 100526
$ws.Range("F14").Value = "Bronshjon"
$ws.Range("G14").Value = "Callidium coriaceum"
$ws.Range("H14").Value = "Paykull, 1800"
$ws.Range("M14").Value = "äldre gnagspår"
$ws.Range("Q14").Value = 703965.55072247
$ws.Range("R14").Value = 6572785.445717536
$ws.Range("AO14").Value = "torrgran"
